# Updated symbol list on Sat Dec 17 07:48:35 UTC 2022 with GitHub Actions
#
# Refresh the cryptocurrency price/volume snapshot values on Sheet1.
# Price ("D") column cells are stored as text (they already contain
# values like "5.300" / "0.1041" with significant trailing zeros), so
# numeric-looking replacements are written with a leading apostrophe to
# keep Excel from coercing them into floating-point numbers. Plain text
# columns (B/C/E) are written as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
}

function Set-PlainValue($addr, $value) {
    $ws.Range($addr).Value = $value
}

# -- Price column refreshes (numeric-looking text, needs quote prefix) --
Set-TextValue "D2"  "235.35"
Set-TextValue "D3"  "22.54"
Set-TextValue "D4"  "5.299"
Set-TextValue "D5"  "0.05609"
Set-TextValue "D6"  "3.378"
Set-TextValue "D7"  "6.483"
Set-TextValue "D8"  "1.067"
Set-TextValue "D9"  "0.7811"
Set-TextValue "D10" "0.1395"
Set-TextValue "D11" "0.07392"
Set-TextValue "D12" "0.03167"
Set-TextValue "D13" "0.02969"
Set-TextValue "D14" "0.09269"
Set-TextValue "D15" "0.001659"
Set-TextValue "D16" "3.252"

Set-TextValue "D18" "0.0005799"
Set-PlainValue "E18" "17OneONEWorstin24h"

Set-TextValue "D19" "0.006235"
Set-TextValue "D20" "0.005227"
Set-TextValue "D21" "0.001050"
Set-TextValue "D22" "0.0001498"
Set-TextValue "D23" "3.972"

Set-TextValue "D26" "0.1055"
Set-TextValue "D27" "0.0004989"

Set-TextValue "D40" "0.04053"
Set-TextValue "D41" "0.006941"

# -- Rows 42/43: BKEXToken and CEJI swap positions --
Set-PlainValue "B42" "CEJI"
Set-PlainValue "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue  "D42" "0.003497"
Set-PlainValue "E42" "41CEJICEJI"

Set-PlainValue "B43" "BKEXToken"
Set-PlainValue "C43" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue  "D43" "0.1040"
Set-PlainValue "E43" "42BKEXTokenBKK"

Set-TextValue "D44" "0.01034"
Set-TextValue "D45" "0.00005439"
Set-TextValue "D46" "0.00000000749"
Set-TextValue "D47" "0.6751"

Set-TextValue "D48" "0.04058"
Set-PlainValue "E48" "47BOLOBOLO"

Set-TextValue "D49" "0.00002098"
Set-TextValue "D50" "0.01009"
